$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.089.82"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "1.651.42"
$ws.Range("E3").Value = "  -1.04%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "'218.54"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'0.5246"
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  -0.51%  "
$ws.Range("D8").Value = "'0.2674"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").Value = "'0.07685"
$ws.Range("E11").Value = "  -1.93%  "
$ws.Range("D12").Value = "'4.588"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("D13").Value = "1.662.63"
$ws.Range("E13").Value = "  -0.33%  "
$ws.Range("D14").Value = "1.879.22"
$ws.Range("E14").Value = "  -0.93%  "
$ws.Range("D15").Value = "'0.5617"
$ws.Range("E15").Value = "  -0.11%  "
$ws.Range("D16").Value = "0.0₅8223"
$ws.Range("D17").Value = "'65.43"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").Value = "26.091.68"
$ws.Range("E18").Value = "  -0.90%  "
$ws.Range("E19").Value = "  -0.53%  "
$ws.Range("D20").Value = "'4.678"
$ws.Range("E20").Value = "  -0.84%  "
$ws.Range("B21").Value = "BitcoinCash"
$ws.Range("C21").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D21").Value = "'190.88"
$ws.Range("E21").Value = "  -4.37%  "
$ws.Range("B22").Value = "Avalanche"
$ws.Range("C22").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D22").Value = "'10.33"
$ws.Range("E22").Value = "  +0.32%  "
$ws.Range("D23").Value = "'5.970"
$ws.Range("E23").Value = "  -1.45%  "
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("D25").Value = "'146.13"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -1.16%  "
$ws.Range("D27").Value = "'7.233"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").Value = "'15.94"
$ws.Range("E28").Value = "  -1.35%  "
$ws.Range("D29").Value = "'1.497"
$ws.Range("E29").Value = "  -1.73%  "
$ws.Range("D30").Value = "'0.05625"
$ws.Range("E30").Value = "  -4.49%  "
$ws.Range("D31").Value = "'1.270"
$ws.Range("E31").Value = "  -1.19%  "
$ws.Range("D32").Value = "'3.490"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "'3.380"
$ws.Range("E33").Value = "  +1.61%  "
$ws.Range("D34").Value = "'1.574"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("D35").Value = "'2.791"
$ws.Range("E35").Value = "  -1.41%  "
$ws.Range("D36").Value = "'2.410"
$ws.Range("E36").Value = "  -0.92%  "
$ws.Range("D37").Value = "'0.9442"
$ws.Range("E37").Value = "  -2.00%  "
$ws.Range("D38").Value = "'0.5768"
$ws.Range("E38").Value = "  -0.44%  "
$ws.Range("E39").Value = "  -1.56%  "
$ws.Range("D40").Value = "'5.966"
$ws.Range("E40").Value = "  -0.05%  "
$ws.Range("E41").Value = "  -0.55%  "
$ws.Range("D42").Value = "'0.8395"
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "1.020.04"
$ws.Range("E43").Value = "  -5.39%  "
$ws.Range("D44").Value = "'101.35"
$ws.Range("E44").Value = "  -1.54%  "
$ws.Range("D45").Value = "1.790.75"
$ws.Range("E45").Value = "  -0.89%  "
$ws.Range("D46").Value = "'58.43"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05354"
$ws.Range("E47").Value = "  +3.97%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₈104"
$ws.Range("E48").Value = "  +0.71%  "
$ws.Range("E49").Value = "  -1.28%  "
$ws.Range("D50").Value = "'8.043"
$ws.Range("E50").Value = "  -0.49%  "
$ws.Range("D51").Value = "'0.4342"
$ws.Range("E51").Value = "  -1.72%  "
